# Add rows 2-5 with session dates (stored as literal text, not real dates).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store these date-look-alike strings as
# plain text instead of auto-converting them to date serial numbers.
$ws.Range("A2").Value = "'2025-08-28"
$ws.Range("A3").Value = "'2025-08-26"
$ws.Range("A4").Value = "'2025-08-27"
$ws.Range("A5").Value = "'2025-08-25"

# Reset to the default "Normal" style so the new cells pick up no explicit
# formatting (the apostrophe trick otherwise marks them with a quote-prefix
# style).
$ws.Range("A2:A5").Style = "Normal"
